$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNum, D(Fecha), I(Calidad), J(Volumen), K(PrecioMin), L(PrecioMax), M(PrecioPromPond), N(Unidad), P(PrecioKg), Q(KgOUnidades)
$rows = @(
    @(2, 44536, 'Primera', 790, 14000, 15000, 14494, '$/caja 13 kilos', 1115, 13),
    @(3, 44536, 'Segunda', 430, 11000, 12000, 11500, '$/caja 13 kilos', 885, 13),
    @(6, 44396, 'Primera', 770, 17000, 18000, 17494, '$/caja 13 kilos', 1346, 13),
    @(7, 44396, 'Segunda', 340, 16000, 16000, 16000, '$/caja 13 kilos', 1231, 13),
    @(8, 44326, 'Primera', 340, 25000, 25000, 25000, '$/caja 13 kilos', 1923, 13),
    @(9, 44326, 'Segunda', 160, 23000, 23000, 23000, '$/caja 13 kilos', 1769, 13),
    @(10, 44270, 'Primera', 250, 40000, 40000, 40000, '$/caja 15 kilos', 2667, 15),
    @(11, 44382, 'Primera', 790, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
    @(12, 44382, 'Segunda', 430, 12000, 12000, 12000, '$/caja 13 kilos', 923, 13),
    @(13, 44459, 'Primera', 970, 13000, 14000, 13495, '$/caja 13 kilos', 1038, 13),
    @(14, 44459, 'Segunda', 520, 11000, 12000, 11500, '$/caja 13 kilos', 885, 13),
    @(15, 44235, 'Primera', 250, 42000, 43000, 42400, '$/caja 13 kilos', 3262, 13),
    @(16, 44249, 'Primera', 250, 39000, 42000, 40500, '$/caja 13 kilos', 3115, 13),
    @(17, 44165, 'Primera', 430, 31000, 32000, 31465, '$/caja 13 kilos', 2420, 13),
    @(18, 44403, 'Primera', 700, 16000, 17000, 16500, '$/caja 13 kilos', 1269, 13),
    @(19, 44403, 'Segunda', 430, 15000, 15000, 15000, '$/caja 13 kilos', 1154, 13),
    @(20, 44557, 'Primera', 970, 17000, 18000, 17495, '$/caja 13 kilos', 1346, 13),
    @(21, 44557, 'Segunda', 430, 16000, 16000, 16000, '$/caja 13 kilos', 1231, 13),
    @(22, 44550, 'Primera', 790, 11000, 12000, 11494, '$/caja 13 kilos', 884, 13),
    @(23, 44550, 'Segunda', 430, 9000, 10000, 9500, '$/caja 13 kilos', 731, 13),
    @(24, 44410, 'Primera', 790, 15000, 16000, 15494, '$/caja 13 kilos', 1192, 13),
    @(25, 44410, 'Segunda', 340, 13000, 13000, 13000, '$/caja 13 kilos', 1000, 13),
    @(26, 44487, 'Primera', 1150, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
    @(27, 44487, 'Segunda', 610, 12000, 12000, 12000, '$/caja 13 kilos', 923, 13),
    @(28, 44333, 'Primera', 340, 25000, 26000, 25500, '$/caja 13 kilos', 1962, 13),
    @(29, 44333, 'Segunda', 160, 23000, 23000, 23000, '$/caja 13 kilos', 1769, 13),
    @(30, 44340, 'Primera', 250, 20000, 20000, 20000, '$/caja 13 kilos', 1538, 13),
    @(31, 44340, 'Segunda', 160, 18000, 18000, 18000, '$/caja 13 kilos', 1385, 13),
    @(32, 44291, 'Primera', 340, 24000, 25000, 24500, '$/caja 13 kilos', 1885, 13),
    @(33, 44371, 'Primera', 160, 20000, 21000, 20500, '$/caja 13 kilos', 1577, 13),
    @(34, 44263, 'Primera', 250, 40000, 40000, 40000, '$/caja 15 kilos', 2667, 15),
    @(35, 44200, 'Primera', 520, 30000, 30000, 30000, '$/caja 13 kilos', 2308, 13),
    @(36, 44200, 'Segunda', 340, 25000, 25000, 25000, '$/caja 13 kilos', 1923, 13),
    @(37, 44435, 'Primera', 880, 13000, 14000, 13500, '$/caja 13 kilos', 1038, 13),
    @(38, 44435, 'Segunda', 340, 11000, 12000, 11500, '$/caja 13 kilos', 885, 13),
    @(39, 44284, 'Primera', 340, 28000, 30000, 29000, '$/caja 13 kilos', 2231, 13),
    @(40, 44242, 'Primera', 250, 41000, 43000, 42000, '$/caja 13 kilos', 3231, 13),
    @(41, 44424, 'Primera', 700, 13000, 14000, 13500, '$/caja 13 kilos', 1038, 13),
    @(42, 44424, 'Segunda', 430, 12000, 12000, 12000, '$/caja 13 kilos', 923, 13),
    @(43, 44508, 'Primera', 1150, 15000, 16000, 15500, '$/caja 13 kilos', 1192, 13),
    @(44, 44417, 'Primera', 790, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
    @(45, 44417, 'Segunda', 340, 13000, 13000, 13000, '$/caja 13 kilos', 1000, 13),
    @(46, 44529, 'Primera', 790, 16000, 18000, 16987, '$/caja 13 kilos', 1307, 13),
    @(47, 44529, 'Segunda', 430, 13000, 14000, 13500, '$/caja 13 kilos', 1038, 13),
    @(48, 44172, 'Primera', 430, 30000, 30000, 30000, '$/caja 13 kilos', 2308, 13),
    @(49, 44319, 'Primera', 340, 24000, 25000, 24500, '$/caja 13 kilos', 1885, 13),
    @(50, 44319, 'Segunda', 160, 22000, 22000, 22000, '$/caja 13 kilos', 1692, 13),
    @(51, 44179, 'Primera', 430, 29000, 30000, 29465, '$/caja 13 kilos', 2267, 13),
    @(52, 44438, 'Primera', 790, 13000, 14000, 13494, '$/caja 13 kilos', 1038, 13),
    @(53, 44438, 'Segunda', 340, 11000, 12000, 11500, '$/caja 13 kilos', 885, 13),
    @(54, 44431, 'Primera', 880, 13000, 14000, 13500, '$/caja 13 kilos', 1038, 13),
    @(55, 44431, 'Segunda', 340, 11000, 12000, 11500, '$/caja 13 kilos', 885, 13),
    @(56, 44186, 'Primera', 450, 29000, 30000, 29444, '$/caja 13 kilos', 2265, 13),
    @(57, 44389, 'Primera', 700, 19000, 20000, 19500, '$/caja 13 kilos', 1500, 13),
    @(58, 44389, 'Segunda', 340, 17000, 17000, 17000, '$/caja 13 kilos', 1308, 13),
    @(59, 44221, 'Primera', 350, 40000, 42000, 40857, '$/caja 13 kilos', 3143, 13),
    @(60, 44221, 'Segunda', 180, 35000, 35000, 35000, '$/caja 13 kilos', 2692, 13),
    @(61, 44515, 'Primera', 1060, 16000, 18000, 17000, '$/caja 13 kilos', 1308, 13),
    @(62, 44515, 'Segunda', 610, 14000, 14000, 14000, '$/caja 13 kilos', 1077, 13),
    @(63, 44473, 'Primera', 1060, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
    @(64, 44473, 'Segunda', 430, 11000, 12000, 11500, '$/caja 13 kilos', 885, 13),
    @(65, 44298, 'Primera', 340, 24000, 25000, 24500, '$/caja 13 kilos', 1885, 13),
    @(66, 44305, 'Primera', 340, 24000, 24000, 24000, '$/caja 13 kilos', 1846, 13),
    @(67, 44305, 'Segunda', 160, 20000, 20000, 20000, '$/caja 13 kilos', 1538, 13),
    @(68, 44445, 'Primera', 790, 13000, 14000, 13494, '$/caja 13 kilos', 1038, 13),
    @(69, 44445, 'Segunda', 340, 11000, 12000, 11500, '$/caja 13 kilos', 885, 13),
    @(70, 44466, 'Primera', 1150, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
    @(71, 44466, 'Segunda', 790, 12000, 12000, 12000, '$/caja 13 kilos', 923, 13),
    @(72, 44494, 'Primera', 780, 15000, 15000, 15000, '$/caja 13 kilos', 1154, 13),
    @(73, 44571, 'Primera', 610, 12000, 13000, 12500, '$/caja 13 kilos', 962, 13),
    @(74, 44571, 'Segunda', 106, 10000, 10000, 10000, '$/caja 13 kilos', 769, 13),
    @(75, 44277, 'Primera', 250, 38000, 38000, 38000, '$/caja 13 kilos', 2923, 13),
    @(76, 44277, 'Segunda', 160, 35000, 35000, 35000, '$/caja 13 kilos', 2692, 13),
    @(77, 44354, 'Primera', 250, 15000, 16000, 15500, '$/caja 13 kilos', 1192, 13),
    @(78, 44312, 'Primera', 430, 25000, 25000, 25000, '$/caja 13 kilos', 1923, 13),
    @(79, 44312, 'Segunda', 250, 23000, 23000, 23000, '$/caja 13 kilos', 1769, 13)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D Fecha
    $ws.Cells.Item($r, 9).Value  = $row[2]   # I Calidad
    $ws.Cells.Item($r, 10).Value = $row[3]   # J Volumen
    $ws.Cells.Item($r, 11).Value = $row[4]   # K Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[5]   # L Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[6]   # M Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $row[7]   # N Unidad de comercializacion
    $ws.Cells.Item($r, 16).Value = $row[8]   # P Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $row[9]   # Q Kg o Unidades
}
